$d = $word.ActiveDocument

# MOIS DU Juilet 2020 -> MOIS DU Aout 2020
$d.Content.Find.Execute("Juilet", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Aout", 2)

# Code / reference number 4159 -> 4190
$d.Content.Find.Execute("4159", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4190", 2)

# Total base amount 74 862 000,00 -> 75 420 000,00
$d.Content.Find.Execute("74 862 000,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "75 420 000,00", 2)

# Cotisation amounts (appears twice) 3 743 100,00 -> 3 771 000,00
$d.Content.Find.Execute("3 743 100,00", $true, $false, $false, $false, $false,
                         $true, 2, $false, "3 771 000,00", 2)

# Spelled-out amount in words
$d.Content.Find.Execute("TROIS MILLIONS SEPT CENT QUARANTE-TROIS MILLE CENT ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TROIS MILLIONS SEPT CENT SOIXANTE ET ONZE MILLE  ", 2)
